# Applies cell value updates per the diff for Jogos_da_Semana_FlashScore_2025-02-21.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 1.83
$ws.Range("I8").Value = 4.75
$ws.Range("J8").Value = 2.5
$ws.Range("L8").Value = 5.5
$ws.Range("N8").Value = 8
$ws.Range("W8").Value = 1.5
$ws.Range("X8").Value = 2.5
$ws.Range("AL8").Value = 10

# Row 9
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 2.55
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 4
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = 1.2
$ws.Range("N9").Value = 4.33
$ws.Range("O9").Value = 1.91
$ws.Range("P9").Value = 1.8
$ws.Range("Q9").Value = 4.2
$ws.Range("R9").Value = 1.22
$ws.Range("S9").Value = 8
$ws.Range("T9").Value = 1.08
$ws.Range("U9").Value = 10
$ws.Range("V9").Value = 1.06
$ws.Range("W9").Value = 1.93
$ws.Range("X9").Value = 1.88
$ws.Range("Y9").Value = 2.75
$ws.Range("Z9").Value = 1.4
$ws.Range("AA9").Value = 5.5
$ws.Range("AB9").Value = 12
$ws.Range("AD9").Value = 34
$ws.Range("AG9").Value = 4
$ws.Range("AJ9").Value = 151
$ws.Range("AL9").Value = 5.5
$ws.Range("AM9").Value = 12
$ws.Range("AN9").Value = 15
$ws.Range("AO9").Value = 34
$ws.Range("AP9").Value = 41
$ws.Range("AQ9").Value = 67
$ws.Range("AR9").Value = 3.05
$ws.Range("AS9").Value = 1.39

# Row 10
$ws.Range("G10").Value = 2.2
$ws.Range("I10").Value = 3.6
$ws.Range("J10").Value = 3.1
$ws.Range("L10").Value = 4.33
$ws.Range("S10").Value = 4.4
$ws.Range("W10").Value = 1.62
$ws.Range("X10").Value = 2.2
$ws.Range("AB10").Value = 9
$ws.Range("AD10").Value = 21
$ws.Range("AE10").Value = 23
$ws.Range("AG10").Value = 6
$ws.Range("AL10").Value = 8
$ws.Range("AN10").Value = 13
$ws.Range("AP10").Value = 34

# Row 25
$ws.Range("G25").Value = 1.91
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 4.75
$ws.Range("J25").Value = 2.75
$ws.Range("L25").Value = 5
$ws.Range("M25").Value = 1.11
$ws.Range("N25").Value = 6.5
$ws.Range("S25").Value = 4.1
$ws.Range("U25").Value = 5
$ws.Range("V25").Value = 1.17
$ws.Range("Y25").Value = 2.2
$ws.Range("Z25").Value = 1.62
$ws.Range("AA25").Value = 5.5
$ws.Range("AB25").Value = 7.5
$ws.Range("AD25").Value = 15
$ws.Range("AH25").Value = 6
$ws.Range("AL25").Value = 9.5
$ws.Range("AM25").Value = 21
$ws.Range("AN25").Value = 17
$ws.Range("AP25").Value = 41
$ws.Range("AQ25").Value = 51
$ws.Range("AR25").Value = 1.95
$ws.Range("AS25").Value = 1.9

# Row 31
$ws.Range("G31").Value = 2.35
$ws.Range("H31").Value = 2.9
$ws.Range("I31").Value = 3.25
$ws.Range("J31").Value = 3.1
$ws.Range("K31").Value = 2.05
$ws.Range("M31").Value = 1.07
$ws.Range("N31").Value = 9
$ws.Range("O31").Value = 1.36
$ws.Range("P31").Value = 3
$ws.Range("Q31").Value = 2.15
$ws.Range("R31").Value = 1.67
$ws.Range("U31").Value = 3.75
$ws.Range("V31").Value = 1.25
$ws.Range("AD31").Value = 21
$ws.Range("AG31").Value = 8

# Row 60
$ws.Range("G60").Value = 1.62
$ws.Range("H60").Value = 4.1
$ws.Range("K60").Value = 2.4
$ws.Range("Q60").Value = 1.65
$ws.Range("R60").Value = 2.2
$ws.Range("Y60").Value = 1.67
$ws.Range("AL60").Value = 17
$ws.Range("AM60").Value = 29
$ws.Range("AN60").Value = 17
$ws.Range("AP60").Value = 41
$ws.Range("AQ60").Value = 41

# Row 70
$ws.Range("G70").Value = 2.5
$ws.Range("I70").Value = 3
$ws.Range("J70").Value = 3.4
$ws.Range("Q70").Value = 2.5
$ws.Range("R70").Value = 1.5
$ws.Range("U70").Value = 5
$ws.Range("V70").Value = 1.17
$ws.Range("W70").Value = 1.57
$ws.Range("X70").Value = 2.25
$ws.Range("Y70").Value = 2.1
$ws.Range("Z70").Value = 1.67
$ws.Range("AC70").Value = 11
$ws.Range("AD70").Value = 26
$ws.Range("AG70").Value = 6.5
$ws.Range("AK70").Value = 1000
$ws.Range("AL70").Value = 7.5
$ws.Range("AM70").Value = 13
$ws.Range("AO70").Value = 29
$ws.Range("AR70").Value = 1.85
$ws.Range("AS70").Value = 1.95

# Row 71
$ws.Range("G71").Value = 1.2
$ws.Range("H71").Value = 5.5
$ws.Range("I71").Value = 15
$ws.Range("K71").Value = 2.6
$ws.Range("Q71").Value = 1.65
$ws.Range("R71").Value = 2.2
$ws.Range("W71").Value = 1.3
$ws.Range("X71").Value = 3.4
$ws.Range("AC71").Value = 9
$ws.Range("AF71").Value = 29
$ws.Range("AG71").Value = 13
$ws.Range("AH71").Value = 11
$ws.Range("AI71").Value = 26
$ws.Range("AN71").Value = 29
$ws.Range("AR71").Value = 1.3
$ws.Range("AS71").Value = 3.6
